# Scenario_Component_Vehicle.xlsx edit:
# Add a second vehicle data row (row 3) to the OperationScenario_Vehicle sheet,
# expand Table1 to cover the new row (plus one trailing blank table row, as
# Excel does when a row is appended right after a table), and move the
# active cell selection to D8 to match the saved state from the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 3) - mirrors row 2's layout/units, with zeroed-out
# numeric fields for the new vehicle entry.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "electricity"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "Wh"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "Wh/km"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "W"
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "W"
$ws.Range("M3").Value = 0

# Expand the table (Table1) to include the new row plus one extra blank
# row, matching the resulting table ref of A1:M4 while sheetData only goes
# to row 3.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:M4"))

# Update the selected cell shown when the workbook was last saved.
$ws.Range("D8").Select()
